$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    83.703703703703695,
    83.703703703703695,
    83.703703703703695,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    83.333333333333343,
    83.333333333333343,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    83.333333333333343,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    83.333333333333343,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481,
    81.481481481481481
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}
